$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "97.139.35"

$ws.Range("E2").Value = "  +4.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.128.70"

$ws.Range("E3").Value = "  +0.59%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.25"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  +2.52%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "612.00"
$ws.Range("D6").ClearFormats()

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.10"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = "  +2.06%  "

# Row 8 - Dogecoin
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.384"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "  -1.28%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.127.02"

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.34%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.02%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "96.785.05"

$ws.Range("E13").Value = "  +4.62%  "

# Row 14 - ShibaInu
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000240"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = "  -1.89%  "

# Row 15 - Avalanche
$ws.Range("B15").Value = "Toncoin"

$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.46"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = "  +0.63%  "

# Row 16 - Toncoin
$ws.Range("B16").Value = "Avalanche"

$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.03"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  +0.12%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.711.42"

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.117.14"

$ws.Range("E18").Value = "  +0.22%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "526.50"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  +19.87%  "

# Row 20 - SuiNetwork
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.51"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  -7.76%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.56"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = "  +0.67%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  -2.27%  "

# Row 23 - PEPE
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000192"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -5.53%  "

# Row 24 - Uniswap
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.83"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  -3.08%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.57"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = "  +3.66%  "

# Row 26 - NEARProtocol
$ws.Range("E26").Value = "  -2.12%  "

# Row 27 - Aptos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.61"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = "  -1.42%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "3.297.50"

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.17%  "

# Row 30 - Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.237"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = "  +2.58%  "

# Row 31 - Cronos
$ws.Range("E31").Value = "  -4.73%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -1.64%  "

# Row 33 - Binance-PegBSC-USD
$ws.Range("E33").Value = "  -3.54%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "8.98"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = "  -1.88%  "

# Row 35 - EthereumClassic
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.62"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = "  +3.62%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -5.87%  "

# Row 37 - RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.37"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = "  -9.07%  "

# Row 38 - PancakeSwap
$ws.Range("E38").Value = "  -0.98%  "

# Row 39 - Bittensor
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "482.20"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  +3.41%  "

# Row 40 - WhiteBITCoin
$ws.Range("E40").Value = "  +1.33%  "

# Row 41 - PolygonEcosystemToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.439"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = "  +2.71%  "

# Row 42 - Fetch.AI
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.22"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  -5.03%  "

# Row 43 - MantraDAO
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.58"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  -10.55%  "

# Row 45 - dogwifhat
$ws.Range("E45").Value = "  -5.26%  "

# Row 46 - Monero
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "160.75"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = "  +1.06%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  +4.52%  "

# Row 48 - ARBITRUM
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.699"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  +2.53%  "

# Row 49 - Filecoin
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.46"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = "  +2.83%  "

# Row 50 - OKB
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.43"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = "  +0.75%  "

# Row 51 - FirstDigitalUSD
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = "  +0.04%  "
